$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Body" (paragraph 3) -> text becomes "Lead", moves to list level 2
#    (ilvl=1 in OOXML terms).
# ------------------------------------------------------------------
$pLead = $d.Paragraphs.Item(3)
$rLead = $pLead.Range
$rLead.MoveEnd(1, -1)          # exclude the paragraph mark
$rLead.Text = "Lead"
$pLead.Range.ListFormat.ListLevelNumber = 2

# ------------------------------------------------------------------
# 2. Insert a brand new paragraph right after "Lead" for "Ending"
#    (also at list level 2 / ilvl=1). It will carry the relocated
#    "_GoBack" bookmark, anchored right at the end of its text.
# ------------------------------------------------------------------
$pLead.Range.InsertParagraphAfter()
$pEnding = $d.Paragraphs.Item(4)
$rEnding = $pEnding.Range
$rEnding.MoveEnd(1, -1)
# Type a temporary trailing marker character so the bookmark can be
# anchored at a "safe" (non paragraph-end) position first, then
# deleted -- directly creating a zero-width range right before the
# paragraph mark snaps to the wrong location.
$rEnding.Text = "EndingZ"

# Remove the existing (hidden) "_GoBack" bookmark wherever it is now.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-create "_GoBack" right after "Ending" (before the temporary "Z").
$pEndingStart = $pEnding.Range.Start
$bmPos = $pEndingStart + 6
$rBm = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $rBm)

# Delete the temporary trailing "Z" -- the zero-width bookmark we just
# created stays put because it sits strictly before the deleted text.
$rZ = $d.Range($bmPos, $bmPos + 1)
$rZ.Delete()

# ------------------------------------------------------------------
# 3. The paragraph that used to read "Conclusion" (now shifted down to
#    position 5 because of the new "Ending" paragraph) becomes "Body"
#    (list level 1 / ilvl=0). Its bookmark already moved away above.
# ------------------------------------------------------------------
$pBody2 = $d.Paragraphs.Item(5)
$rBody2 = $pBody2.Range
$rBody2.MoveEnd(1, -1)
$rBody2.Text = "Body"
$pBody2.Range.ListFormat.ListLevelNumber = 1

# ------------------------------------------------------------------
# 4. Insert a new paragraph after it for "Conclusion" (list level 1 /
#    ilvl=0). This paragraph does NOT carry any bookmark.
# ------------------------------------------------------------------
$pBody2.Range.InsertParagraphAfter()
$pConclusion2 = $d.Paragraphs.Item(6)
$rConclusion2 = $pConclusion2.Range
$rConclusion2.MoveEnd(1, -1)
$rConclusion2.Text = "Conclusion"
$pConclusion2.Range.ListFormat.ListLevelNumber = 1

Write-Output "done"
